# Auto-generated Word COM-interop script to apply the two-digit-mul.docx edit
$d = $word.ActiveDocument

# 1) Update the date/weekday heading at the top of the document
$d.Content.Find.Execute("2023-04-05 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-04-06 Thursday", 2) | Out-Null

# 2) Update each "a×b=" multiplication prompt in the 20x5 table, in row-major
#    (row, then column) order -- this mirrors the document order of the diff.
$newValues = @(
    "46×83=",
    "84×56=",
    "82×97=",
    "86×48=",
    "72×57=",
    "67×37=",
    "35×43=",
    "71×95=",
    "25×98=",
    "53×22=",
    "53×28=",
    "75×64=",
    "24×13=",
    "93×10=",
    "34×30=",
    "98×62=",
    "60×41=",
    "96×56=",
    "34×40=",
    "66×65=",
    "31×63=",
    "35×33=",
    "51×77=",
    "61×23=",
    "78×10=",
    "53×78=",
    "44×73=",
    "44×17=",
    "32×75=",
    "86×94=",
    "21×94=",
    "82×92=",
    "74×84=",
    "54×78=",
    "78×95=",
    "40×11=",
    "37×89=",
    "53×60=",
    "15×35=",
    "74×10=",
    "27×29=",
    "95×26=",
    "46×12=",
    "18×36=",
    "91×34=",
    "79×100=",
    "18×26=",
    "14×85=",
    "78×43=",
    "82×39=",
    "19×43=",
    "35×24=",
    "32×44=",
    "64×49=",
    "37×82=",
    "68×42=",
    "38×79=",
    "25×92=",
    "64×88=",
    "15×75=",
    "94×29=",
    "98×28=",
    "30×98=",
    "72×26=",
    "38×52=",
    "48×76=",
    "29×46=",
    "55×63=",
    "13×15=",
    "98×42=",
    "95×74=",
    "38×28=",
    "52×11=",
    "47×69=",
    "60×33=",
    "95×42=",
    "13×78=",
    "27×98=",
    "100×36=",
    "88×78=",
    "73×22=",
    "78×51=",
    "13×84=",
    "21×68=",
    "79×45=",
    "17×29=",
    "71×60=",
    "20×30=",
    "24×19=",
    "36×98=",
    "27×95=",
    "95×31=",
    "52×83=",
    "15×97=",
    "63×97=",
    "95×45=",
    "58×75=",
    "45×68=",
    "70×95=",
    "94×56="
)

$table = $d.Tables.Item(1)
$rows = $table.Rows.Count
$cols = $table.Columns.Count

if (($rows * $cols) -ne $newValues.Length) {
    throw "Expected $($newValues.Length) cells but table has $rows x $cols = $($rows * $cols)"
}

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $table.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Updated date heading and $idx table cells."
